# Fruta / hortaliza, semanal
# Insert a new weekly record at row 147 (Feria Lagunitas de Puerto Montt - Pomelo),
# pushing the existing historical rows (147-167) down by one (to 148-168).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 147; this shifts rows 147:167 down
# to 148:168 (values + formatting), matching the growth of the used range to
# A1:T168.
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new weekly observation.
# (Single-quoted literals throughout so no value is ever treated as a
# PowerShell expandable string — safe for the embedded '$' and apostrophe.)
$ws.Cells.Item(147, 1).Value = 4
$ws.Cells.Item(147, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(147, 3).Value = 'Los Lagos'
$ws.Cells.Item(147, 4).Value = 44505
$ws.Cells.Item(147, 5).Value = 10
$ws.Cells.Item(147, 6).Value = 'Fruta'
$ws.Cells.Item(147, 7).Value = 100102
$ws.Cells.Item(147, 8).Value = 'Cítricos'
$ws.Cells.Item(147, 9).Value = 100102006
$ws.Cells.Item(147, 10).Value = 'Pomelo'
$ws.Cells.Item(147, 11).Value = 'Start Ruby'
$ws.Cells.Item(147, 12).Value = 'Primera'
$ws.Cells.Item(147, 13).Value = 180
$ws.Cells.Item(147, 14).Value = 11000
$ws.Cells.Item(147, 15).Value = 12000
$ws.Cells.Item(147, 16).Value = 11500
$ws.Cells.Item(147, 17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(147, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(147, 19).Value = 821
$ws.Cells.Item(147, 20).Value = 14
